$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 21:20"
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 237877
$ws.Range("C4").Value = 22874
$ws.Range("D4").Value = 10324
$ws.Range("E4").Value = 221835
$ws.Range("F4").Value = 5421
$ws.Range("G4").Value = 616
$ws.Range("H4").Value = 5718
$ws.Range("A7").Value = "Alemania"
$ws.Range("B7").Value = 84600
$ws.Range("C7").Value = 6619
$ws.Range("D7").Value = 21400
$ws.Range("E7").Value = 62103
$ws.Range("F7").Value = 3936
$ws.Range("G7").Value = 166
$ws.Range("H7").Value = 1097
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 11118
$ws.Range("C16").Value = 1387
$ws.Range("D16").Value = 1906
$ws.Range("E16").Value = 9078
$ws.Range("F16").Value = 120
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 134
$ws.Range("A17").Value = "Austria"
$ws.Range("B17").Value = 11108
$ws.Range("C17").Value = 397
$ws.Range("D17").Value = 1749
$ws.Range("E17").Value = 9201
$ws.Range("F17").Value = 227
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 158
$ws.Range("A21").Value = "Israel"
$ws.Range("B21").Value = 6857
$ws.Range("C21").Value = 765
$ws.Range("D21").Value = 338
$ws.Range("E21").Value = 6485
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 8
$ws.Range("H21").Value = 34
$ws.Range("A22").Value = "Suecia"
$ws.Range("B22").Value = 5568
$ws.Range("C22").Value = 621
$ws.Range("D22").Value = 103
$ws.Range("E22").Value = 5157
$ws.Range("F22").Value = 429
$ws.Range("G22").Value = 69
$ws.Range("H22").Value = 308
$ws.Range("A24").Value = "Noruega"
$ws.Range("B24").Value = 5136
$ws.Range("C24").Value = 259
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 5054
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 50
$ws.Range("A45").Value = "Peru"
$ws.Range("B45").Value = 1414
$ws.Range("C45").Value = 91
$ws.Range("D45").Value = 537
$ws.Range("E45").Value = 830
$ws.Range("F45").Value = 49
$ws.Range("G45").Value = 9
$ws.Range("H45").Value = 47
$ws.Range("A54").Value = "Emiratos Arabes Unidos"
$ws.Range("B54").Value = 1024
$ws.Range("C54").Value = 210
$ws.Range("D54").Value = 96
$ws.Range("E54").Value = 920
$ws.Range("F54").Value = 2
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 8
$ws.Range("A58").Value = "Ucrania"
$ws.Range("B58").Value = 897
$ws.Range("C58").Value = 103
$ws.Range("D58").Value = 19
$ws.Range("E58").Value = 856
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 22
$ws.Range("A59").Value = "Eslovenia"
$ws.Range("B59").Value = 897
$ws.Range("C59").Value = 56
$ws.Range("D59").Value = 70
$ws.Range("E59").Value = 810
$ws.Range("F59").Value = 31
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 17
$ws.Range("A60").Value = "Estonia"
$ws.Range("B60").Value = 858
$ws.Range("C60").Value = 79
$ws.Range("D60").Value = 45
$ws.Range("E60").Value = 802
$ws.Range("F60").Value = 16
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 11
$ws.Range("A61").Value = "Egipto"
$ws.Range("B61").Value = 850
$ws.Range("C61").Value = 71
$ws.Range("D61").Value = 179
$ws.Range("E61").Value = 619
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 52
$ws.Range("A81").Value = "Costa Rica"
$ws.Range("B81").Value = 396
$ws.Range("C81").Value = 21
$ws.Range("D81").Value = 6
$ws.Range("E81").Value = 388
$ws.Range("F81").Value = 11
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 2
$ws.Range("A82").Value = "Republica de Macedonia"
$ws.Range("B82").Value = 384
$ws.Range("C82").Value = 30
$ws.Range("D82").Value = 17
$ws.Range("E82").Value = 356
$ws.Range("F82").Value = 4
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 11
$ws.Range("A93").Value = "Afganistan"
$ws.Range("B93").Value = 273
$ws.Range("C93").Value = 36
$ws.Range("D93").Value = 10
$ws.Range("E93").Value = 257
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 6
$ws.Range("A104").Value = "Nigeria"
$ws.Range("B104").Value = 184
$ws.Range("C104").Value = 10
$ws.Range("D104").Value = 20
$ws.Range("E104").Value = 162
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 2
$ws.Range("A105").Value = "Islas Feroe"
$ws.Range("B105").Value = 177
$ws.Range("C105").Value = 4
$ws.Range("D105").Value = 81
$ws.Range("E105").Value = 96
$ws.Range("F105").Value = 1
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("A111").Value = "Martinica"
$ws.Range("B111").Value = 138
$ws.Range("C111").Value = 3
$ws.Range("D111").Value = 27
$ws.Range("E111").Value = 108
$ws.Range("F111").Value = 19
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 3
$ws.Range("A190").Value = "San Bartolome"
$ws.Range("B190").Value = 6
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 1
$ws.Range("E190").Value = 5
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0
$ws.Range("A191").Value = "Cabo Verde"
$ws.Range("B191").Value = 6
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 5
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 1
$ws.Range("A192").Value = "Nepal"
$ws.Range("B192").Value = 6
$ws.Range("C192").Value = 1
$ws.Range("D192").Value = 1
$ws.Range("E192").Value = 5
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0
$ws.Range("A193").Value = "Mauritania"
$ws.Range("B193").Value = 6
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 2
$ws.Range("E193").Value = 3
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 1
$ws.Range("A194").Value = "Islas Turcas y Caicos"
$ws.Range("B194").Value = 5
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 0
$ws.Range("E194").Value = 5
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0
